$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in cell A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Agosto de 2020 a las 13:36"

# Apply country ranking re-sort and updated case statistics
$ws.Cells.Item(4, 5).Value = 2271253
$ws.Cells.Item(4, 7).Value = 7
$ws.Cells.Item(4, 8).Value = 161608
$ws.Cells.Item(14, 2).Value = 320117
$ws.Cells.Item(14, 3).Value = 2634
$ws.Cells.Item(14, 4).Value = 277463
$ws.Cells.Item(14, 5).Value = 24678
$ws.Cells.Item(14, 7).Value = 174
$ws.Cells.Item(14, 8).Value = 17976
$ws.Cells.Item(68, 1).Value = "Nepal"
$ws.Cells.Item(68, 2).Value = 21750
$ws.Cells.Item(68, 3).Value = 360
$ws.Cells.Item(68, 4).Value = 15389
$ws.Cells.Item(68, 5).Value = 6296
$ws.Cells.Item(68, 7).Value = 5
$ws.Cells.Item(68, 8).Value = 65
$ws.Cells.Item(69, 1).Value = "Austria"
$ws.Cells.Item(69, 2).Value = 21696
$ws.Cells.Item(69, 3).Value = 130
$ws.Cells.Item(69, 4).Value = 19596
$ws.Cells.Item(69, 5).Value = 1381
$ws.Cells.Item(69, 8).Value = 719
$ws.Cells.Item(79, 1).Value = "Estado de Palestina"
$ws.Cells.Item(79, 2).Value = 13398
$ws.Cells.Item(79, 3).Value = 333
$ws.Cells.Item(79, 4).Value = 6907
$ws.Cells.Item(79, 5).Value = 6402
$ws.Cells.Item(79, 8).Value = 89
$ws.Cells.Item(80, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(80, 2).Value = 13306
$ws.Cells.Item(80, 3).Value = 168
$ws.Cells.Item(80, 4).Value = 7031
$ws.Cells.Item(80, 5).Value = 5891
$ws.Cells.Item(80, 7).Value = 5
$ws.Cells.Item(80, 8).Value = 384
$ws.Cells.Item(82, 2).Value = 12526
$ws.Cells.Item(82, 3).Value = 304
$ws.Cells.Item(82, 4).Value = 10148
$ws.Cells.Item(82, 5).Value = 2244
$ws.Cells.Item(82, 7).Value = 7
$ws.Cells.Item(82, 8).Value = 134
$ws.Cells.Item(85, 2).Value = 10715
$ws.Cells.Item(85, 3).Value = 177
$ws.Cells.Item(85, 4).Value = 7101
$ws.Cells.Item(85, 5).Value = 3391
$ws.Cells.Item(85, 7).Value = 5
$ws.Cells.Item(85, 8).Value = 223
$ws.Cells.Item(126, 5).Value = 189
$ws.Cells.Item(126, 7).Value = 1
$ws.Cells.Item(126, 8).Value = 125
$ws.Cells.Item(128, 2).Value = 2124
$ws.Cells.Item(128, 3).Value = 11
$ws.Cells.Item(128, 4).Value = 1954
$ws.Cells.Item(128, 5).Value = 107
$ws.Cells.Item(134, 2).Value = 1930
$ws.Cells.Item(134, 3).Value = 4
$ws.Cells.Item(134, 5).Value = 95
$ws.Cells.Item(150, 1).Value = "Malta"
$ws.Cells.Item(150, 2).Value = 946
$ws.Cells.Item(150, 3).Value = 20
$ws.Cells.Item(150, 4).Value = 670
$ws.Cells.Item(150, 5).Value = 267
$ws.Cells.Item(150, 8).Value = 9
$ws.Cells.Item(151, 2).Value = 944
$ws.Cells.Item(151, 4).Value = 296
$ws.Cells.Item(151, 5).Value = 600
$ws.Cells.Item(151, 8).Value = 48
$ws.Cells.Item(152, 4).Value = 825
$ws.Cells.Item(152, 5).Value = 62
$ws.Cells.Item(152, 8).Value = 52
$ws.Cells.Item(153, 2).Value = 939
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(153, 4).Value = 835
$ws.Cells.Item(153, 5).Value = 29
$ws.Cells.Item(153, 8).Value = 75
$ws.Cells.Item(154, 1).Value = "Siria"
$ws.Cells.Item(154, 2).Value = 928
$ws.Cells.Item(154, 3).Value = 8
$ws.Cells.Item(154, 4).Value = 745
$ws.Cells.Item(154, 5).Value = 171
$ws.Cells.Item(154, 8).Value = 12
$ws.Cells.Item(159, 1).Value = "Vietnam"
$ws.Cells.Item(159, 2).Value = 747
$ws.Cells.Item(159, 3).Value = 30
$ws.Cells.Item(159, 4).Value = 392
$ws.Cells.Item(159, 5).Value = 345
$ws.Cells.Item(159, 8).Value = 10
$ws.Cells.Item(160, 1).Value = "Lesoto"
$ws.Cells.Item(160, 2).Value = 742
$ws.Cells.Item(160, 3).Value = 16
$ws.Cells.Item(160, 4).Value = 175
$ws.Cells.Item(160, 5).Value = 544
$ws.Cells.Item(160, 8).Value = 23
$ws.Cells.Item(180, 2).Value = 190
$ws.Cells.Item(180, 3).Value = 1
$ws.Cells.Item(180, 4).Value = 184
$ws.Cells.Item(202, 1).Value = "Santa Lucia"
$ws.Cells.Item(203, 1).Value = "Timor Oriental"
